$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded ahead of the existing
# entries, so insert a fresh row at 176 (this shifts every row from
# 176..265 down to 177..266, preserving their original values/styles).
$ws.Rows(176).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A176").Value = 10
$ws.Range("B176").Value = "Vega Modelo de Temuco"
$ws.Range("C176").Value = "La Araucanía"
$ws.Range("D176").Value = 45001
$ws.Range("E176").Value = 9
$ws.Range("F176").Value = 100114007
$ws.Range("G176").Value = "Jengibre"
$ws.Range("H176").Value = "Sin especificar"
$ws.Range("I176").Value = "Primera"
$ws.Range("J176").Value = 150
$ws.Range("K176").Value = 25000
$ws.Range("L176").Value = 25000
$ws.Range("M176").Value = 25000
$ws.Range("N176").Value = "$/caja 13 kilos"
$ws.Range("O176").Value = "Perú"
$ws.Range("P176").Value = 1923
$ws.Range("Q176").Value = 13
$ws.Range("R176").Value = "Hortaliza"
